$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.428.93"
$ws.Range("D3").Value = "1.840.52"
$ws.Range("E3").Value = "  -2.09%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'260.17"
$ws.Range("E5").Value = "  -6.68%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "'0.5213"
$ws.Range("E7").Value = "  -1.77%  "
$ws.Range("D8").Value = "'0.3240"
$ws.Range("E8").Value = "  -6.42%  "
$ws.Range("D9").Value = "'0.06779"
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("D10").Value = "'18.55"
$ws.Range("E10").Value = "  -8.10%  "
$ws.Range("D11").Value = "'0.7641"
$ws.Range("E11").Value = "  -5.51%  "
$ws.Range("D12").Value = "'0.07672"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("D13").Value = "1.865.61"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").Value = "'88.24"
$ws.Range("E14").Value = "  -2.84%  "
$ws.Range("D15").Value = "'5.017"
$ws.Range("E15").Value = "  -3.36%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "'13.90"
$ws.Range("E17").Value = "  -4.83%  "
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "'0.000007946"
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("D20").Value = "26.460.71"
$ws.Range("D21").Value = "2.073.84"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("D22").Value = "'4.557"
$ws.Range("E22").Value = "  -4.19%  "
$ws.Range("D23").Value = "'9.457"
$ws.Range("E23").Value = "  -6.17%  "
$ws.Range("E24").Value = "  -4.21%  "
$ws.Range("D25").Value = "'144.20"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").Value = "'2.196"
$ws.Range("E26").Value = "  -7.31%  "
$ws.Range("D27").Value = "'1.657"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").Value = "'16.94"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("D29").Value = "'110.97"
$ws.Range("E29").Value = "  -2.38%  "
$ws.Range("D30").Value = "'4.152"
$ws.Range("E30").Value = "  -5.18%  "
$ws.Range("D31").Value = "'4.117"
$ws.Range("E31").Value = "  -5.08%  "
$ws.Range("D32").Value = "'0.08721"
$ws.Range("E32").Value = "  -2.03%  "
$ws.Range("D33").Value = "'0.04782"
$ws.Range("E33").Value = "  -3.42%  "
$ws.Range("E34").Value = "  -5.09%  "
$ws.Range("D35").Value = "'2.849"
$ws.Range("E35").Value = "  -1.28%  "
$ws.Range("D36").Value = "'0.6967"
$ws.Range("E36").Value = "  -5.06%  "
$ws.Range("D37").Value = "'3.064"
$ws.Range("E37").Value = "  -7.03%  "
$ws.Range("D38").Value = "'0.01748"
$ws.Range("E38").Value = "  -5.94%  "
$ws.Range("D39").Value = "'2.176"
$ws.Range("E39").Value = "  -9.00%  "
$ws.Range("D40").Value = "'0.4807"
$ws.Range("E40").Value = "  -7.00%  "
$ws.Range("D41").Value = "'110.71"
$ws.Range("E41").Value = "  -4.82%  "
$ws.Range("D42").Value = "'0.8903"
$ws.Range("E42").Value = "  -7.15%  "
$ws.Range("D43").Value = "'6.053"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("D44").Value = "'1.002"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").Value = "'7.629"
$ws.Range("E45").Value = "  -6.16%  "
$ws.Range("D46").Value = "'0.4120"
$ws.Range("E46").Value = "  -8.91%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05855"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.962"
$ws.Range("E48").Value = "  -4.26%  "
$ws.Range("D50").Value = "'0.1216"
$ws.Range("E50").Value = "  -9.84%  "
$ws.Range("D51").Value = "'0.8815"
$ws.Range("E51").Value = "  -0.53%  "
